$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.315.56"
$ws.Range("E2").Value = "  +3.55%  "

# Row 3
$ws.Range("D3").Value = "3.488.80"
$ws.Range("E3").Value = "  +2.96%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.59%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.65%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "3.490.46"
$ws.Range("E9").Value = "  +3.04%  "

# Row 10
$ws.Range("E10").Value = "  -1.10%  "

# Row 11
$ws.Range("E11").Value = "  +3.59%  "

# Row 12
$ws.Range("E12").Value = "  +3.61%  "

# Row 13
$ws.Range("D13").Value = "4.091.48"
$ws.Range("E13").Value = "  +3.24%  "

# Row 14
$ws.Range("E14").Value = "  +1.38%  "

# Row 15
$ws.Range("E15").Value = "  +3.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.70%  "

# Row 17
$ws.Range("D17").Value = "65.343.48"
$ws.Range("E17").Value = "  +3.59%  "

# Row 18
$ws.Range("D18").Value = "3.488.12"
$ws.Range("E18").Value = "  +3.06%  "

# Row 19
$ws.Range("E19").Value = "  +3.81%  "

# Row 20
$ws.Range("E20").Value = "  +2.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.94%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.554"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "

# Row 26
$ws.Range("E26").Value = "  +6.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.03%  "

# Row 28
$ws.Range("E28").Value = "  +2.25%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.67%  "

# Row 32
$ws.Range("E32").Value = "  +4.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

# Row 34
$ws.Range("E34").Value = "  +8.80%  "

# Row 35
$ws.Range("E35").Value = "  +10.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.35%  "

# Row 37
$ws.Range("E37").Value = "  +6.43%  "

# Row 38
$ws.Range("D38").Value = "3.014.88"
$ws.Range("E38").Value = "  +3.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0779"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.45%  "

# Row 41
$ws.Range("E41").Value = "  +2.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.781"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.94%  "

# Row 47
$ws.Range("E47").Value = "  +5.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "324.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.69%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.37%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.28%  "
